# =====================================================================
# PlayerPerformance_3699.xlsx restructuring
#   - Insert a new "Player Info" sheet at the front
#   - Rename MATCH_CARD_LINK -> MATCH_CODE on "ODI Batting"/"ODI Bowling"
#     and replace the howstat URL values with the bare match code
#   - Drop the (always-empty) INNING_NUMBER cells on "ODI Batting" for
#     matches the player did not bat in
#   - Append a new "ODI Batting Extra" sheet with additional per-match
#     batting detail
# =====================================================================

$wb = $excel.ActiveWorkbook

function Set-HeaderCell($cell, $text) {
    $cell.Value = $text
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

function Set-TextCell($cell, $text) {
    # Force text storage so numeric-looking strings ("0", "14.22%", ...)
    # are not silently reinterpreted as numbers / percentages.
    $cell.NumberFormat = "@"
    $cell.Value = [string]$text
}

# ---------------------------------------------------------------------
# 1. Rename MATCH_CARD_LINK -> MATCH_CODE and collapse the howstat URL
#    down to the bare match code, on both existing sheets.
# ---------------------------------------------------------------------

$batting = $wb.Worksheets.Item("ODI Batting")
$bowling = $wb.Worksheets.Item("ODI Bowling")

# ODI Batting: MATCH_CARD_LINK lives in column D
$battingLinkCol = 4
Set-HeaderCell ($batting.Cells.Item(1, $battingLinkCol)) "MATCH_CODE"

$battingRowCount = $batting.UsedRange.Rows.Count
for ($r = 2; $r -le $battingRowCount; $r++) {
    $cell = $batting.Cells.Item($r, $battingLinkCol)
    $url = $cell.Text
    if ($url -and $url -ne "") {
        $code = ($url -split "MatchCode=")[1]
        Set-TextCell $cell $code
    }
}

# ODI Bowling: MATCH_CARD_LINK lives in column B
$bowlingLinkCol = 2
Set-HeaderCell ($bowling.Cells.Item(1, $bowlingLinkCol)) "MATCH_CODE"

$bowlingRowCount = $bowling.UsedRange.Rows.Count
for ($r = 2; $r -le $bowlingRowCount; $r++) {
    $cell = $bowling.Cells.Item($r, $bowlingLinkCol)
    $url = $cell.Text
    if ($url -and $url -ne "") {
        $code = ($url -split "MatchCode=")[1]
        Set-TextCell $cell $code
    }
}

# ---------------------------------------------------------------------
# 2. ODI Batting: the INNING_NUMBER column (B) is blank whenever the
#    player did not bat in the match - those cells should not exist at
#    all (rather than hold an empty string).
# ---------------------------------------------------------------------

for ($r = 2; $r -le $battingRowCount; $r++) {
    $inningCell = $batting.Cells.Item($r, 2)
    $dismissal = $batting.Cells.Item($r, 8).Text
    if ($dismissal -eq "did not bat") {
        $inningCell.ClearContents()
    }
}

# ---------------------------------------------------------------------
# 3. Add the "Player Info" sheet at the very front of the workbook.
# ---------------------------------------------------------------------

$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"
$playerInfo.Move($batting)

Set-HeaderCell ($playerInfo.Cells.Item(1, 1)) "ID"
Set-HeaderCell ($playerInfo.Cells.Item(1, 2)) "NAME"
Set-HeaderCell ($playerInfo.Cells.Item(1, 3)) "BATTING_HAND"
Set-HeaderCell ($playerInfo.Cells.Item(1, 4)) "BOWL_STYLE"

Set-TextCell ($playerInfo.Cells.Item(2, 1)) "3699"
Set-TextCell ($playerInfo.Cells.Item(2, 2)) "Isuru Udana Tillakaratna"
Set-TextCell ($playerInfo.Cells.Item(2, 3)) "Right Handed"
Set-TextCell ($playerInfo.Cells.Item(2, 4)) "Left Arm Fast Medium"

# ---------------------------------------------------------------------
# 4. Add the "ODI Batting Extra" sheet at the very end of the workbook.
# ---------------------------------------------------------------------

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$battingExtra = $wb.Worksheets.Add($null, $lastSheet)
$battingExtra.Name = "ODI Batting Extra"

Set-HeaderCell ($battingExtra.Cells.Item(1, 1)) "MATCH_CODE"
Set-HeaderCell ($battingExtra.Cells.Item(1, 2)) "BATTING_POSITION"
Set-HeaderCell ($battingExtra.Cells.Item(1, 3)) "NUM_4"
Set-HeaderCell ($battingExtra.Cells.Item(1, 4)) "NUM_6"
Set-HeaderCell ($battingExtra.Cells.Item(1, 5)) "PERCENT_RUNS_OF_TOTAL"
Set-HeaderCell ($battingExtra.Cells.Item(1, 6)) "MAN_OF_MATCH"

$extraRows = @(
    @("3433", "10", "",  "",  "",       "NO"),
    @("4269", "7",  "0", "0", "",       "NO"),
    @("4271", "",   "",  "",  "",       "NO"),
    @("4272", "9",  "3", "1", "14.22%", "NO"),
    @("4302", "",   "",  "",  "",       "NO"),
    @("4305", "9",  "0", "0", "",       "NO"),
    @("4309", "8",  "0", "1", "4.98%",  "NO"),
    @("4322", "9",  "1", "0", "3.24%",  "NO"),
    @("4331", "9",  "1", "0", "2.59%",  "NO"),
    @("4339", "9",  "1", "0", "8.37%",  "NO"),
    @("4344", "7",  "0", "0", "0.89%",  "NO"),
    @("4350", "9",  "0", "0", "0.38%",  "NO"),
    @("4357", "",   "",  "",  "",       "NO"),
    @("4375", "9",  "0", "0", "0.42%",  "NO"),
    @("4413", "9",  "0", "0", "",       "NO"),
    @("4414", "",   "",  "",  "",       "NO"),
    @("4417", "9",  "0", "0", "0.65%",  "NO"),
    @("4463", "9",  "2", "0", "9.38%",  "NO"),
    @("4464", "",   "",  "",  "",       "NO"),
    @("4480", "9",  "0", "0", "3.05%",  "NO")
)

$r = 2
foreach ($row in $extraRows) {
    Set-TextCell ($battingExtra.Cells.Item($r, 1)) $row[0]

    $posCell = $battingExtra.Cells.Item($r, 2)
    if ($row[1] -eq "") {
        $posCell.NumberFormat = "@"
        $posCell.Value = ""
    } else {
        $posCell.Value = [int]$row[1]
    }

    Set-TextCell ($battingExtra.Cells.Item($r, 3)) $row[2]
    Set-TextCell ($battingExtra.Cells.Item($r, 4)) $row[3]
    Set-TextCell ($battingExtra.Cells.Item($r, 5)) $row[4]
    Set-TextCell ($battingExtra.Cells.Item($r, 6)) $row[5]

    $r++
}

Write-Host "Done restructuring workbook."
